# Update Allegro Group office address in the Ralph CLA document.
$d = $word.ActiveDocument

# Replace the street name "Marcelińska" with "Grunwaldzka".
$d.Content.Find.Execute("Marcelińska", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Grunwaldzka", 2)

# Replace the old house number / postal code / city fragment with the new one.
$d.Content.Find.Execute(" 90, 60-324 ", $true, $false, $false, $false, $false,
                         $true, 1, $false, " 182, 60-166 ", 2)
